$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (INSTRUCTIONS): German text tightened; row height shrinks as a
#     consequence of the shorter wrapped text (150 -> 135 pts). ---
$newInstructionsDe = "Im Folgenden werden Ihnen Sätze präsentiert, die **einfache Sachverhalte** beschreiben. Der Inhalt der Sätze kann wahr oder falsch sein. Bitte entscheiden Sie **so korrekt** aber auch **so schnell** wie möglich, ob es sich um eine wahre oder um eine falsche Aussage handelt.`n"
$ws.Range("B3").Value = $newInstructionsDe
$ws.Range("B3").EntireRow.RowHeight = 135

# --- Row 4 (INSTRUCTIONS2): German text gains parenthetical explanations
#     for the J / F keys. ---
$newInstructions2De = "Bei **wahren** Aussagen drücken Sie bitte die **“J”-Taste** (`"J`" für `"ja, wahr`").\\`n Bei **falschen** Aussagen drücken Sie bitte die **“F”-Taste** (`"F`" für `"falsch`") .\\**Bitte legen Sie Ihre Zeigefinger jetzt auf die “F”- und “J”-Taste**.\\ "
$ws.Range("B4").Value = $newInstructions2De

# --- Row 10 (WELCOME): German text replaced. ---
$ws.Range("B10").Value = "Test: Leseverständnis"

# --- Cursor/selection left on A17 when the file was saved. ---
$ws.Range("A17").Select()
